$wb = $excel.ActiveWorkbook

# --- "Constant Samples" sheet: move the selection cursor from B9 to B6 ---
$wsConstant = $wb.Worksheets.Item("Constant Samples")
$wsConstant.Range("B6").Select() | Out-Null

# --- "Formula Samples_CheckOrder" sheet: add the "prereq" cell D4 = "k" ---
# (written first so that "k" claims shared-string slot 53, before "h" is
# introduced on the other sheet, matching the string table ordering)
$wsFormulaCheckOrder = $wb.Worksheets.Item("Formula Samples_CheckOrder")
$wsFormulaCheckOrder.Range("D4").Value = "k"
$wsFormulaCheckOrder.Range("D4").Select() | Out-Null

# --- "Constant Samples_CheckOrder" sheet: add the killer/prereq column D ---
$wsConstantCheckOrder = $wb.Worksheets.Item("Constant Samples_CheckOrder")
$wsConstantCheckOrder.Range("D4").Value = "h"
$wsConstantCheckOrder.Range("D5").Value = "k"
$wsConstantCheckOrder.Range("D6").Value = "k"

# This sheet becomes the active / selected tab, with D7 as the active cell
$wsConstantCheckOrder.Activate() | Out-Null
$wsConstantCheckOrder.Range("D7").Select() | Out-Null
